$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the original segment names (currently in column A, rows 2-20)
# before we shift columns around.
$segments = @()
for ($r = 2; $r -le 20; $r++) {
    $segments += $ws.Cells.Item($r, 1).Value()
}

# Insert a new column before column B. This shifts the old B:K headers /
# data to C:L, and the old A (segment names) stays where it is for now.
$ws.Columns.Item(2).Insert()

# Build the new "segments" header in B1, copying the bold/centered/bordered
# header formatting from the (now shifted) C1 header cell.
$ws.Range("C1").Copy()
$ws.Range("B1").PasteSpecial(-4122)
$ws.Range("B1").Value = "segments"

# Column A becomes a numeric index (0-based), keeping its existing style.
for ($r = 2; $r -le 20; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 2
}

# Column B gets the segment names as plain (unstyled) text.
$ws.Range("B2:B20").ClearFormats()
for ($r = 2; $r -le 20; $r++) {
    $ws.Cells.Item($r, 2).Value = $segments[$r - 2]
}
